$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.535.55'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.823.78'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.89'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5180'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3925'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07689'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.00'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.112'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.00'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.79%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.279'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("B14").Value = 'BinanceUSD'
$ws.Range("C14").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.004'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.541'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.824.97'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.58'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001079'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06606'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.68'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.62%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.059'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.547.46'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.13'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.243'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +7.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.60'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.27'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.035.98'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.418'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.91'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.133'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1108'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.654'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.656'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07239'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2243'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.98%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.952'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.04%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02333'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.150'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.30'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6246'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.182'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.37'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5906'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.85%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.715'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.65'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.188'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06936'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.75%  '
